$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 42; $r++) {
    # Column A ("Score"): round the numeric score to 2 decimal places
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    $aCell.Value = [Math]::Round($aVal, 2)

    # Column K ("Rent Budget"): convert "(570, 988)" style text into "£570-£988"
    $kCell = $ws.Cells.Item($r, 11)
    $kVal = $kCell.Value2
    if ($kVal -match '\((\d+),\s*(\d+)\)') {
        $low = $matches[1]
        $high = $matches[2]
        $kCell.Value = "£" + $low + "-£" + $high
    }
}
